# Apply the commit's edits to assets_liabilities.xlsx
#   - Summary sheet: update name, income, totals, net worth, ratio
#   - Assets sheet:  insert "Vehicles" rows (Luxury Car, Economy Car),
#                    update Liquid Assets + TOTAL ASSETS values
#   - Liabilities:   insert "Auto Loans" row (Vehicle Loan 1),
#                    update Personal Loan / Credit Card / TOTAL values

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Faisal Al Ameri"
$summary.Range("B4").Value = 3584.14
$summary.Range("B6").Value = 625184
$summary.Range("B7").Value = 463394
$summary.Range("B8").Value = 161790
$summary.Range("B9").Value = 1.35

# ---------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

# Insert two new rows above the existing "Liquid Assets" row (row 2),
# pushing "Liquid Assets" to row 4 and "TOTAL ASSETS" to row 5.
$assets.Range("2:3").Insert()

$newRows = $assets.Range("A2:C3")
$newRows.ClearFormats()
$newRows.Borders.LineStyle = 1
$assets.Range("C2:C3").NumberFormat = "#,##0"

$assets.Range("A2").Value = "Vehicles"
$assets.Range("B2").Value = "Luxury Car"
$assets.Range("C2").Value = 569032

$assets.Range("A3").Value = "Vehicles"
$assets.Range("B3").Value = "Economy Car"
$assets.Range("C3").Value = 49845

# Row 4 is the former row 2 ("Liquid Assets" / "Savings Account") - value changed
$assets.Range("C4").Value = 6307

# Row 5 is the former row 3 ("TOTAL ASSETS") - value changed
$assets.Range("C5").Value = 625184

# ---------------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------------
$liabilities = $wb.Worksheets.Item("Liabilities")

# Insert one new row above the existing "Personal Loans" row (row 2),
# pushing "Personal Loans" to row 3, "Credit Cards" to row 4 and
# "TOTAL LIABILITIES" to row 5.
$liabilities.Range("2:2").Insert()

$newRow = $liabilities.Range("A2:E2")
$newRow.ClearFormats()
$newRow.Borders.LineStyle = 1
$liabilities.Range("C2:D2").NumberFormat = "#,##0"

$liabilities.Range("A2").Value = "Auto Loans"
$liabilities.Range("B2").Value = "Vehicle Loan 1"
$liabilities.Range("C2").Value = 341419
$liabilities.Range("D2").Value = 7113
$liabilities.Range("E2").Value = 4

# Row 3 is the former row 2 ("Personal Loans" / "Personal Loan") - values changed
$liabilities.Range("C3").Value = 94635
$liabilities.Range("D3").Value = 1972
$liabilities.Range("E3").Value = 4

# Row 4 is the former row 3 ("Credit Cards" / "Credit Card Balance") - values changed
$liabilities.Range("C4").Value = 27340
$liabilities.Range("D4").Value = 1367
$liabilities.Range("E4").Value = 1

# Row 5 is the former row 4 ("TOTAL LIABILITIES") - value changed
$liabilities.Range("C5").Value = 463394
